# AutoCommit_25 февраля 2024 г. 16:21:30_SibNout2023
#
# This script reproduces the edit described by the diff:
#  - Header row (row 2, C2:G2) switches from date columns to homework
#    columns "ДЗ_1".."ДЗ_5"; H2 loses its header text entirely.
#  - All "Н" marks and stray "5" values inside the attendance/grade grid
#    (rows 4-31, columns C-H) are cleared back to blank (style kept).
#  - The totals row (row 32) is cleared back to blank (style kept).
#  - The frozen pane view is scrolled down so row 16 becomes the first
#    visible (non-frozen) row, and the header row C2:G2 becomes selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row: replace the six date headers with five "ДЗ_n" headers
#    and leave H2 blank.
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "ДЗ_1"
$ws.Range("D2").Value = "ДЗ_2"
$ws.Range("E2").Value = "ДЗ_3"
$ws.Range("F2").Value = "ДЗ_4"
$ws.Range("G2").Value = "ДЗ_5"
$ws.Range("H2").Value = ""

# ---------------------------------------------------------------------
# 2. Clear every "Н" mark in the grade grid (rows 4-31), and the stray
#    "5" values left in the grid. (Each cell is addressed individually
#    -- multi-area Range(...).Value assignment only reliably touches
#    the first area in this runtime, so we avoid relying on that.)
# ---------------------------------------------------------------------
$clearCells = "D5","C7","C9","D9","E11","C12","D12","C15","C21","C25","D25","D27","E27","E4","E23","E28"
foreach ($addr in $clearCells) {
    $ws.Range($addr).Value = ""
}

# ---------------------------------------------------------------------
# 4. Clear the totals row (row 32).
# ---------------------------------------------------------------------
$ws.Range("C32:H32").Value = ""

# ---------------------------------------------------------------------
# 5. Update the frozen-pane view: keep the header rows (1-3) and the
#    first two columns (A-B) frozen, but scroll the data pane so that
#    row 16 is the first visible row, then select the new header range.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("C16").Select()
$win.FreezePanes = $true
$ws.Range("C2:G2").Select()
